$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (tab) from "economiedata" to "data"
$ws.Name = "data"

# Remove the (now-unused) external reference link from the workbook
$wb.BreakLink("/office/statistiek-excel/resources/grafieken.xlsx", 1)

# Update the active selection to A28
$ws.Range("A28").Select()
